# Update run metadata to reflect a new pipeline execution:
#  - "preguntas" sheet: refresh the per-row timestamp (column T, rows 2-23)
#  - "indice_global" sheet: refresh run_id (L2), timestamp (Q2), and the
#    run_id embedded inside the extra_config dict string (R2)

$wb = $excel.ActiveWorkbook

$oldTimestamp = "2025-10-26T11:24:53.518499"
$newTimestamp = "2025-10-26T11:29:42.909132"

$oldRunId = "32a45db9b48e454bb75cb54551a0162d"
$newRunId = "f8e202da056c46b48a6e0608d49a82d5"

# --- preguntas: update timestamp column (T) for every data row ---
$preguntas = $wb.Worksheets.Item("preguntas")
for ($row = 2; $row -le 23; $row++) {
    $cell = $preguntas.Cells.Item($row, 20)  # column T = 20
    if ([string]$cell.Value2 -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}

# --- indice_global: update run_id, timestamp, and extra_config text ---
$indice = $wb.Worksheets.Item("indice_global")

$indice.Range("L2").Value = $newRunId
$indice.Range("Q2").Value = $newTimestamp

$extraConfig = [string]$indice.Range("R2").Value2
$extraConfig = $extraConfig -replace [regex]::Escape($oldRunId), $newRunId
$indice.Range("R2").Value = $extraConfig
